$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Currency number format used by column B already; extend it to column E ---
$ws.Range("E2:E20").NumberFormat = "#,##0.00 ""€"""

# --- Move the "nPM1300-QEAA-R7 total" row out of row 16 down to the Total row (20) ---
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()

$ws.Range("D20").Value = "Total"
$ws.Range("E20").Formula = "=SUM(E2:E14)"

# --- New "Delta" columns next to the Total row ---
$ws.Range("G20").Value = "Delta"
$ws.Range("H20").Formula = "=B20-E20"
$ws.Range("I20").Formula = "=E20/B20"
$ws.Range("I20").NumberFormat = "0%"

$ws.Range("A1").Select()
